$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove bold/bordered header style from row 1 (A1:BK1) so cells become unstyled
$ws.Range("A1:BK1").ClearFormats()

# Row 4
$ws.Range("A4").Value = "Phạm Thế Anh"
$ws.Range("B4").Value = "B20DCDT017"
$ws.Range("C4").Value = "D20DTMT1"
$ws.Range("D4").Value = "Nguyễn Trung Hiếu"
$ws.Range("K4").Value = "Nguyễn Quốc Dinh"
$ws.Range("R4").Value = "Nguyễn Quốc Uy"
$ws.Range("Y4").Value = "Lê Thanh Bằng"
$ws.Range("AF4").Value = "Đinh Quang Ngọc"
$ws.Range("AM4").Value = "Nguyễn Trung Hiếu"
$ws.Range("BE4").Value = "Nguyễn Quốc Dinh"

$ws.Range("E4:J4").NumberFormat = "@"
$ws.Range("E4").Value = "3"
$ws.Range("F4").Value = "3"
$ws.Range("G4").Value = "3"
$ws.Range("H4").Value = "3"
$ws.Range("I4").Value = "3"
$ws.Range("J4").Value = "3.0"
$ws.Range("E4:J4").Style = "Normal"
$ws.Range("L4:Q4").NumberFormat = "@"
$ws.Range("L4").Value = "4"
$ws.Range("M4").Value = "4"
$ws.Range("N4").Value = "4"
$ws.Range("O4").Value = "4"
$ws.Range("P4").Value = "4"
$ws.Range("Q4").Value = "4.0"
$ws.Range("L4:Q4").Style = "Normal"
$ws.Range("S4:X4").NumberFormat = "@"
$ws.Range("S4").Value = "5"
$ws.Range("T4").Value = "5"
$ws.Range("U4").Value = "5"
$ws.Range("V4").Value = "5"
$ws.Range("W4").Value = "5"
$ws.Range("X4").Value = "5.0"
$ws.Range("S4:X4").Style = "Normal"
$ws.Range("Z4:AE4").NumberFormat = "@"
$ws.Range("Z4").Value = "6"
$ws.Range("AA4").Value = "6"
$ws.Range("AB4").Value = "6"
$ws.Range("AC4").Value = "6"
$ws.Range("AD4").Value = "6"
$ws.Range("AE4").Value = "6.0"
$ws.Range("Z4:AE4").Style = "Normal"
$ws.Range("AG4:AL4").NumberFormat = "@"
$ws.Range("AG4").Value = "7"
$ws.Range("AH4").Value = "7"
$ws.Range("AI4").Value = "7"
$ws.Range("AJ4").Value = "7"
$ws.Range("AK4").Value = "7"
$ws.Range("AL4").Value = "7.0"
$ws.Range("AG4:AL4").Style = "Normal"
$ws.Range("AN4:AQ4").NumberFormat = "@"
$ws.Range("AN4").Value = "1"
$ws.Range("AO4").Value = "1"
$ws.Range("AP4").Value = "1"
$ws.Range("AQ4").Value = "1.0"
$ws.Range("AN4:AQ4").Style = "Normal"
$ws.Range("BF4:BK4").NumberFormat = "@"
$ws.Range("BF4").Value = "2"
$ws.Range("BG4").Value = "2"
$ws.Range("BH4").Value = "2"
$ws.Range("BI4").Value = "2"
$ws.Range("BJ4").Value = "2"
$ws.Range("BK4").Value = "2.0"
$ws.Range("BF4:BK4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = "Nguyễn Tiến Duy"
$ws.Range("B5").Value = "B20DCDT037"
$ws.Range("C5").Value = "D20DTMT1"
$ws.Range("D5").Value = "Nguyễn Trung Hiếu"
$ws.Range("K5").Value = "Nguyễn Quốc Dinh"
$ws.Range("R5").Value = "Nguyễn Quốc Uy"
$ws.Range("Y5").Value = "Lê Thanh Bằng"
$ws.Range("AF5").Value = "Đinh Quang Ngọc"

$ws.Range("E5:J5").NumberFormat = "@"
$ws.Range("E5").Value = "3"
$ws.Range("F5").Value = "3"
$ws.Range("G5").Value = "3"
$ws.Range("H5").Value = "3"
$ws.Range("I5").Value = "3"
$ws.Range("J5").Value = "3.0"
$ws.Range("E5:J5").Style = "Normal"
$ws.Range("L5:Q5").NumberFormat = "@"
$ws.Range("L5").Value = "4"
$ws.Range("M5").Value = "4"
$ws.Range("N5").Value = "4"
$ws.Range("O5").Value = "4"
$ws.Range("P5").Value = "4"
$ws.Range("Q5").Value = "4.0"
$ws.Range("L5:Q5").Style = "Normal"
$ws.Range("S5:X5").NumberFormat = "@"
$ws.Range("S5").Value = "5"
$ws.Range("T5").Value = "5"
$ws.Range("U5").Value = "5"
$ws.Range("V5").Value = "5"
$ws.Range("W5").Value = "5"
$ws.Range("X5").Value = "5.0"
$ws.Range("S5:X5").Style = "Normal"
$ws.Range("Z5:AE5").NumberFormat = "@"
$ws.Range("Z5").Value = "6"
$ws.Range("AA5").Value = "6"
$ws.Range("AB5").Value = "6"
$ws.Range("AC5").Value = "6"
$ws.Range("AD5").Value = "6"
$ws.Range("AE5").Value = "6.0"
$ws.Range("Z5:AE5").Style = "Normal"
$ws.Range("AG5:AL5").NumberFormat = "@"
$ws.Range("AG5").Value = "7"
$ws.Range("AH5").Value = "7"
$ws.Range("AI5").Value = "7"
$ws.Range("AJ5").Value = "7"
$ws.Range("AK5").Value = "7"
$ws.Range("AL5").Value = "7.0"
$ws.Range("AG5:AL5").Style = "Normal"

# Row 6
$ws.Range("A6").Value = "Lê Sỹ Sang"
$ws.Range("B6").Value = "B20DCDT175"
$ws.Range("C6").Value = "D20DTMT2"
$ws.Range("D6").Value = "Nguyễn Trung Hiếu"
$ws.Range("K6").Value = "Nguyễn Quốc Dinh"
$ws.Range("R6").Value = "Nguyễn Quốc Uy"
$ws.Range("Y6").Value = "Lê Thanh Bằng"
$ws.Range("AF6").Value = "Đinh Quang Ngọc"

$ws.Range("E6:J6").NumberFormat = "@"
$ws.Range("E6").Value = "3"
$ws.Range("F6").Value = "3"
$ws.Range("G6").Value = "3"
$ws.Range("H6").Value = "3"
$ws.Range("I6").Value = "3"
$ws.Range("J6").Value = "3.0"
$ws.Range("E6:J6").Style = "Normal"
$ws.Range("L6:Q6").NumberFormat = "@"
$ws.Range("L6").Value = "4"
$ws.Range("M6").Value = "4"
$ws.Range("N6").Value = "4"
$ws.Range("O6").Value = "4"
$ws.Range("P6").Value = "4"
$ws.Range("Q6").Value = "4.0"
$ws.Range("L6:Q6").Style = "Normal"
$ws.Range("S6:X6").NumberFormat = "@"
$ws.Range("S6").Value = "5"
$ws.Range("T6").Value = "5"
$ws.Range("U6").Value = "5"
$ws.Range("V6").Value = "5"
$ws.Range("W6").Value = "5"
$ws.Range("X6").Value = "5.0"
$ws.Range("S6:X6").Style = "Normal"
$ws.Range("Z6:AE6").NumberFormat = "@"
$ws.Range("Z6").Value = "6"
$ws.Range("AA6").Value = "6"
$ws.Range("AB6").Value = "6"
$ws.Range("AC6").Value = "6"
$ws.Range("AD6").Value = "6"
$ws.Range("AE6").Value = "6.0"
$ws.Range("Z6:AE6").Style = "Normal"
$ws.Range("AG6:AL6").NumberFormat = "@"
$ws.Range("AG6").Value = "7"
$ws.Range("AH6").Value = "7"
$ws.Range("AI6").Value = "7"
$ws.Range("AJ6").Value = "7"
$ws.Range("AK6").Value = "7"
$ws.Range("AL6").Value = "7.0"
$ws.Range("AG6:AL6").Style = "Normal"
